# Generate Report for Handback
# The file "82b901f4-3ba9-4ba9-b5ca-a91d40e19d36.md" has been handed back
# (it was previously "Ready for handoff"). Update its status on every
# sheet that tracks it, and stamp the new handback datetimes for the
# zh-cn and de-de locales.

$wb = $excel.ActiveWorkbook

$overview = $wb.Worksheets.Item("Overview")
$zhcn     = $wb.Worksheets.Item("zh-cn")
$dede     = $wb.Worksheets.Item("de-de")

# Overview sheet: status columns for the zh-cn (B) and de-de (C) locales
$overview.Range("B3").Value = "Handed back: in sync with en-US"
$overview.Range("C3").Value = "Handed back: in sync with en-US"

# zh-cn detail sheet: Status + Latest Handback DateTime
$zhcn.Range("C3").Value = "Handed back: in sync with en-US"
$zhcn.Range("H3").Value = "2016-03-23 10:53:07"

# de-de detail sheet: Status + Latest Handback DateTime
$dede.Range("C3").Value = "Handed back: in sync with en-US"
$dede.Range("H3").Value = "2016-03-23 10:53:13"
